$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) and the "SC 92" row (originally row 28,
# now row 27 after the first deletion) so subsequent rows shift up.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Update column E values that changed between missing/present states.
$ws.Range("E19").Value = -6.5
$ws.Range("E21").Value = ""
$ws.Range("E23").Value = -7
$ws.Range("E27").Value = ""
$ws.Range("E33").Value = -10.7
